$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the document to the
#    empty paragraph right after "... os tópicos:" (paragraph 4).
# ------------------------------------------------------------------

# Remove the existing bookmark at the end of the document first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-add it around the empty paragraph that follows the introductory
# paragraph ("Foram discutidos na reunião online de ontem ... tópicos:").
$paraCount = $d.Paragraphs.Count
$introIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*os tópicos:*") {
        $introIndex = $i
    }
}
if ($introIndex -gt 0) {
    $emptyPara = $d.Paragraphs.Item($introIndex + 1).Range
    $d.Bookmarks.Add("_GoBack", $emptyPara)
}

# ------------------------------------------------------------------
# 2) "Dr" -> "Dra" for the "Shark" character description.
# ------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Shark – Dr ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $drRange1 = $d.Range($rng1.End - 3, $rng1.End - 1)
    $drRange1.Text = "Dra"
}

# ------------------------------------------------------------------
# 3) "Dr" -> "Dra" for the second "Química" bullet point.
# ------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Dr Química - ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $drRange2 = $d.Range($rng2.Start, $rng2.Start + 2)
    $drRange2.Text = "Dra"
}

# ------------------------------------------------------------------
# 4) Merge the two runs that make up "O prazo da segunda entrega é
#    até sábado (10/10/2015)." into a single run/sentence.
# ------------------------------------------------------------------
$d.Content.Find.Execute("O prazo da segunda entrega é até sábado (10/10/2015).", $true, $false, $false, $false, $false, $true, 1, $false, "O prazo da segunda entrega é até sábado (10/10/2015).", 2) | Out-Null
